$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell values scraped from the latest coinranking.com data pull.
$changes = @{
    'D2' = '28.795.63'
    'E2' = '  -0.73%  '
    'D3' = '1.815.57'
    'E3' = '  -0.85%  '
    'D4' = '0.9908'
    'E4' = '  -0.84%  '
    'D5' = '241.59'
    'E5' = '  +0.03%  '
    'D6' = '0.6233'
    'E6' = '  -0.58%  '
    'D7' = '0.9933'
    'E7' = '  -0.73%  '
    'D8' = '0.07400'
    'E8' = '  -2.69%  '
    'D9' = '0.2909'
    'E9' = '  -0.16%  '
    'D10' = '22.83'
    'E10' = '  +0.49%  '
    'D11' = '0.07645'
    'E11' = '  -1.33%  '
    'D12' = '1.813.06'
    'E12' = '  -1.00%  '
    'D13' = '4.961'
    'E13' = '  +0.12%  '
    'D14' = '0.6624'
    'E14' = '  -0.04%  '
    'D15' = '82.41'
    'E15' = '  +0.00%  '
    'D16' = '0.000009517'
    'E16' = '  +0.59%  '
    'D17' = '5.989'
    'E17' = '  +0.10%  '
    'D18' = '28.831.09'
    'E18' = '  -0.55%  '
    'D19' = '12.47'
    'E19' = '  +1.17%  '
    'D20' = '222.76'
    'E20' = '  -0.53%  '
    'D21' = '0.9932'
    'E21' = '  -0.69%  '
    'D22' = '7.058'
    'E22' = '  -2.07%  '
    'D23' = '0.9926'
    'E23' = '  -0.81%  '
    'D24' = '158.55'
    'E24' = '  -0.11%  '
    'D25' = '0.1393'
    'E25' = '  +2.13%  '
    'D26' = '8.432'
    'E26' = '  +0.24%  '
    'D27' = '17.77'
    'E27' = '  -0.32%  '
    'D28' = '1.488'
    'E28' = '  -0.05%  '
    'D29' = '4.087'
    'E29' = '  +0.79%  '
    'D30' = '4.018'
    'E30' = '  -0.13%  '
    'D31' = '0.05417'
    'E31' = '  +4.39%  '
    'D32' = '1.187'
    'E32' = '  -0.40%  '
    'D33' = '1.837'
    'E33' = '  -0.42%  '
    'D34' = '0.7389'
    'E34' = '  +0.20%  '
    'D35' = '1.127'
    'E35' = '  -1.55%  '
    'D36' = '2.585'
    'E36' = '  -4.23%  '
    'D37' = '1.221.51'
    'E37' = '  -2.56%  '
    'D38' = '2.728'
    'E38' = '  -1.20%  '
    'D39' = '0.01768'
    'E39' = '  -1.10%  '
    'D40' = '6.627'
    'E40' = '  +6.70%  '
    'D41' = '0.8921'
    'E41' = '  +0.22%  '
    'D42' = '0.9931'
    'E42' = '  -0.83%  '
    'D43' = '101.11'
    'E43' = '  -0.34%  '
    'D44' = '64.71'
    'E44' = '  +0.25%  '
    'E45' = '  +1.08%  '
    'D46' = '0.5043'
    'E46' = '  -1.28%  '
    'D47' = '0.4010'
    'E47' = '  +0.83%  '
    'D48' = '8.907'
    'E48' = '  +0.81%  '
    'B49' = 'XinFinNetwork'
    'C49' = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
    'D49' = '0.07165'
    'E49' = '  +3.17%  '
    'B50' = 'Cronos'
    'C50' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'D50' = '0.05768'
    'E50' = '  +0.30%  '
    'B51' = 'RenderToken'
    'C51' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D51' = '1.644'
    'E51' = '  +1.23%  '
}

foreach ($cell in $changes.Keys) {
    $range = $ws.Range($cell)
    if ($cell.StartsWith("D")) {
        # Price column: force text so values like "0.07400" or
        # "28.795.63" keep their exact digits/zeros instead of being
        # auto-coerced into a Double (and losing trailing zeros) or,
        # worse, scientific notation.
        $range.NumberFormat = "@"
    }
    $range.Value = $changes[$cell]
}
